# Update YOY data with new column names, preserve original cashflow data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Header renames (Jan-June -> Jan-July) ---
$ws.Range("E1").Value = "Value_2024_Jan_July"
$ws.Range("F1").Value = "Value_2025_Jan_July"

# --- Row 4: VISIT_COUNT ---
$ws.Range("E4").Value = 6263
$ws.Range("F4").Value = 7902
$ws.Range("H4").Value = 26

# --- Row 5: PRIVATE_VISIT_COUNT ---
$ws.Range("E5").Value = 5036
$ws.Range("F5").Value = 6117
$ws.Range("H5").Value = 22

# --- Row 6: STANDARD_COMMITMENTS ---
$ws.Range("E6").Value = 282396
$ws.Range("F6").Value = 215022
$ws.Range("H6").Value = -24

# --- Row 7: VARIABLE_OPERATIONAL_COSTS ---
$ws.Range("E7").Value = 313970
$ws.Range("F7").Value = 305786
$ws.Range("H7").Value = -3

# --- Row 8: MARKETING_ADVERTISING ---
$ws.Range("E8").Value = 56914
$ws.Range("F8").Value = 64665.48
$ws.Range("H8").Value = 14

# --- Row 9: REVENUE_PROPORTIONAL ---
$ws.Range("E9").Value = 70417
$ws.Range("F9").Value = 84980
$ws.Range("H9").Value = 21

# --- Row 10: SUPPLY_PER_VISIT (switch from hardcoded values to formulas) ---
$ws.Range("E10").Formula = "=108025/E4"
$ws.Range("F10").Formula = "=74590/F4"
$ws.Range("H10").Value = -45

# --- Rows 11-13: hardcoded Growth_Rate_Percentage (H) values recompute by hand ---
$ws.Range("H11").Value = -30
$ws.Range("H12").Value = 14
$ws.Range("H13").Value = 1778

# --- View: selection / scroll position ---
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("H4").Select()

$wb.Application.Calculate()
